$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row 770 ("「驚きの価格で」بسعر رائع ..." entry).
# This shifts all subsequent rows (771-794) up by one, matching the diff,
# and reduces the used range from A1:C794 to A1:C793.
$ws.Rows.Item(770).Delete()
